# Fills in rows 14-16 with the new "NMF / LDA" literature entries and
# adds their hyperlinks, then updates the view to match the final state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are assigned in this precise order so that the workbook's shared
# string table is built up in the same sequence as the source edit.
$ws.Range("D14").Value = "https://people.eecs.berkeley.edu/~jfc/hcc/courseSP05/lecs/lec14/NMF03.pdf"
$ws.Range("C14").Value = "Topic Modeling, Aviation Safety, Aviation`nAccident Reports, Machine Learning, LDA, NMF"
$ws.Range("A14").Value = "Topic Modeling Analysis of Aviation Accident Reports: A Comparative Study between LDA and NMF Models"
$ws.Range("B14").Value = "NMF vs LDA"
$ws.Range("E14").Value = 13

$ws.Range("B15").Value = "NMF vs LDA vs BERTopic"
$ws.Range("A15").Value = "Topic Modeling of the SrpELTeC Corpus:`nA Comparison of NMF, LDA, and BERTopic"
$ws.Range("C15").Value = "topic modeling, LDA, NMF, BERTopic, SrpELTeC, computational literary studies"
$ws.Range("D15").Value = "https://ieeexplore.ieee.org/stamp/stamp.jsp?tp=&arnumber=10736093"
$ws.Range("E15").Value = 14

$ws.Range("D16").Value = "https://www.jmlr.org/papers/volume3/blei03a/blei03a.pdf"
$ws.Range("B16").Value = "LDA"
$ws.Range("C16").Value = "LDA"
$ws.Range("A16").Value = "Latent Dirichlet Allocation"
$ws.Range("E16").Value = 15

# Hyperlinks for the new rows
$ws.Hyperlinks.Add($ws.Range("D14"), "https://people.eecs.berkeley.edu/~jfc/hcc/courseSP05/lecs/lec14/NMF03.pdf")
$ws.Hyperlinks.Add($ws.Range("D15"), "https://ieeexplore.ieee.org/stamp/stamp.jsp?tp=&arnumber=10736093")
$ws.Hyperlinks.Add($ws.Range("D16"), "https://www.jmlr.org/papers/volume3/blei03a/blei03a.pdf")

# Update the active view/selection to match the final saved state
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("A17").Select()
